# Add new columns to the LinkML-derived sheets (linkml-sqldb dump todo).
#
# Modification / Modification1: insert a new "parent_parts" column just
# before the existing Modification_id column (which shifts right by one).
#
# Organism / Organism1: insert a new "taxid" column just before the
# existing Organism_id column (which, along with the two columns after it,
# shifts right by one).

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Modification", "Modification1")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns("Q:Q").Insert()
    $ws.Range("Q1").Value = "parent_parts"
}

foreach ($sheetName in @("Organism", "Organism1")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns("E:E").Insert()
    $ws.Range("E1").Value = "taxid"
}
